# edit.ps1
# Applies the row-content rearrangement described by the diff:
# the observation records in rows 2-4, 12-14, and 22-23 were
# re-ordered (their cell values moved between row positions) while
# headers/formatting/other rows stayed untouched.
#
# Row remap (new row <- source/old row):
#   2 <- 4, 3 <- 2, 4 <- 3
#   12 <- 13, 13 <- 14, 14 <- 12
#   22 <- 23, 23 <- 22

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was old row 4)
$ws.Range("A2").Value = 130825852
$ws.Range("P2").Value = "Flinktorpet, Kälom, Offerdal, Jmt"
$ws.Range("Q2").Value = 460952
$ws.Range("R2").Value = 7039723
$ws.Range("S2").Value = 15
$ws.Range("Z2").Value = "10:42"
$ws.Range("AB2").Value = "10:42"
$ws.Range("AC2").Value = "Barkfläkta grövre och klenare granar."

# Row 3 (was old row 2)
$ws.Range("A3").Value = 130826784
$ws.Range("B3").Value = 57884
$ws.Range("E3").Value = 100109
$ws.Range("F3").Value = "Tretåig hackspett"
$ws.Range("G3").Value = "Picoides tridactylus"
$ws.Range("M3").Value = "färska spår"
$ws.Range("P3").Value = "Brännan, Kälom, Offerdal, Jmt"
$ws.Range("Q3").Value = 461233
$ws.Range("R3").Value = 7039438
$ws.Range("Z3").Value = "11:37"
$ws.Range("AB3").Value = "11:37"
$ws.Range("AC3").Value = "Födosök barkfläk"

# Row 4 (was old row 3)
$ws.Range("A4").Value = 130825823
$ws.Range("B4").Value = 57881
$ws.Range("E4").Value = 100049
$ws.Range("F4").Value = "Spillkråka"
$ws.Range("G4").Value = "Dryocopus martius"
$ws.Range("M4").Value = "äldre spår"
$ws.Range("Q4").Value = 460947
$ws.Range("R4").Value = 7039711
$ws.Range("S4").Value = 10
$ws.Range("Z4").Value = "10:38"
$ws.Range("AB4").Value = "10:38"
$ws.Range("AC4").Value = "Födosökshål på äldre döende gran."

# Row 12 (was old row 13)
$ws.Range("A12").Value = 130826287
$ws.Range("B12").Value = 57884
$ws.Range("E12").Value = 100109
$ws.Range("F12").Value = "Tretåig hackspett"
$ws.Range("G12").Value = "Picoides tridactylus"
$ws.Range("H12").Value = "(Linnaeus, 1758)"
$ws.Range("M12").Value = "färska spår"
$ws.Range("P12").Value = "Flinktorpet, Kälom, Offerdal, Jmt"
$ws.Range("Q12").Value = 461096
$ws.Range("R12").Value = 7039690
$ws.Range("Z12").Value = "11:04"
$ws.Range("AB12").Value = "11:04"
$ws.Range("AC12").Value = "Barkfläkta klenare och grövre granar"

# Row 13 (was old row 14)
$ws.Range("A13").Value = 130826478
$ws.Range("P13").Value = "Brännan, Kälom, Offerdal, Jmt"
$ws.Range("Q13").Value = 461220
$ws.Range("R13").Value = 7039590
$ws.Range("Z13").Value = "11:25"
$ws.Range("AB13").Value = "11:25"
$ws.Range("AC13").Value = "Födosök barkfläkt"

# Row 14 (was old row 12)
$ws.Range("A14").Value = 130826137
$ws.Range("B14").Value = 91808
$ws.Range("E14").Value = 1202
$ws.Range("F14").Value = "Ullticka"
$ws.Range("G14").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H14").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("M14").ClearContents()
$ws.Range("P14").Value = "Flinktorpet, Flinktorpet, Jmt"
$ws.Range("Q14").Value = 461026
$ws.Range("R14").Value = 7039757
$ws.Range("Z14").Value = "10:56"
$ws.Range("AB14").Value = "10:56"
$ws.Range("AC14").ClearContents()

# Row 22 (was old row 23)
$ws.Range("A22").Value = 130826438
$ws.Range("B22").Value = 79243
$ws.Range("D22").Value = "NT"
$ws.Range("E22").Value = 6425
$ws.Range("F22").Value = "Garnlav"
$ws.Range("G22").Value = "Alectoria sarmentosa"
$ws.Range("H22").Value = "(Ach.) Ach."
$ws.Range("P22").Value = "Brännan, Brännan, Jmt"
$ws.Range("Q22").Value = 461220
$ws.Range("R22").Value = 7039590
$ws.Range("S22").Value = 25
$ws.Range("Z22").Value = "11:16"
$ws.Range("AB22").Value = "11:16"
$ws.Range("AC22").Value = "Rikligt i området"

# Row 23 (was old row 22)
$ws.Range("A23").Value = 130826355
$ws.Range("B23").Value = 92535
$ws.Range("D23").Value = "VU"
$ws.Range("E23").Value = 67
$ws.Range("F23").Value = "Sprickporing"
$ws.Range("G23").Value = "Diplomitoporus crustulinus"
$ws.Range("H23").Value = "(Bres.) Domański"
$ws.Range("P23").Value = "Flinktorpet, Flinktorpet, Jmt"
$ws.Range("Q23").Value = 461117
$ws.Range("R23").Value = 7039629
$ws.Range("S23").Value = 10
$ws.Range("Z23").Value = "11:10"
$ws.Range("AB23").Value = "11:10"
$ws.Range("AC23").Value = "På undersidan av lutande död gran."
